$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so numeric-looking values
# (e.g. "1.001") are stored as strings, matching the source workbook format.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '30.248.87'
$ws.Range('E2').Value = '  +0.31%  '
$ws.Range('D3').Value = '1.863.24'
$ws.Range('E3').Value = '  +0.01%  '
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').Value = '237.36'
$ws.Range('E5').Value = '  +1.58%  '
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  +0.15%  '
$ws.Range('D7').Value = '0.4680'
$ws.Range('E7').Value = '  +0.37%  '
$ws.Range('D8').Value = '0.2863'
$ws.Range('E8').Value = '  +1.68%  '
$ws.Range('D9').Value = '0.06540'
$ws.Range('E9').Value = '  -0.26%  '
$ws.Range('D10').Value = '22.12'
$ws.Range('E10').Value = '  +11.49%  '
$ws.Range('D11').Value = '0.07886'
$ws.Range('E11').Value = '  +0.92%  '
$ws.Range('D12').Value = '97.66'
$ws.Range('E12').Value = '  +0.99%  '
$ws.Range('D13').Value = '1.866.11'
$ws.Range('E13').Value = '  +0.15%  '
$ws.Range('D14').Value = '5.183'
$ws.Range('E14').Value = '  +0.90%  '
$ws.Range('D15').Value = '0.6796'
$ws.Range('E15').Value = '  +1.79%  '
$ws.Range('D16').Value = '278.24'
$ws.Range('E16').Value = '  -1.11%  '
$ws.Range('D17').Value = '30.258.72'
$ws.Range('E17').Value = '  +0.29%  '
$ws.Range('D18').Value = '13.59'
$ws.Range('E18').Value = '  +7.82%  '
$ws.Range('E19').Value = '  +0.00%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = '0.000007344'
$ws.Range('E20').Value = '  +1.44%  '
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Value = '5.380'
$ws.Range('E21').Value = '  -2.55%  '
$ws.Range('D22').Value = '2.115.06'
$ws.Range('E22').Value = '  +0.71%  '
$ws.Range('E23').Value = '  +0.15%  '
$ws.Range('D24').Value = '6.186'
$ws.Range('E24').Value = '  +0.85%  '
$ws.Range('D25').Value = '168.59'
$ws.Range('E25').Value = '  +1.70%  '
$ws.Range('D26').Value = '9.272'
$ws.Range('E26').Value = '  -0.60%  '
$ws.Range('D27').Value = '19.05'
$ws.Range('E27').Value = '  +0.77%  '
$ws.Range('D28').Value = '1.940'
$ws.Range('E28').Value = '  +1.44%  '
$ws.Range('D29').Value = '1.384'
$ws.Range('E29').Value = '  +3.01%  '
$ws.Range('D30').Value = '0.09799'
$ws.Range('E30').Value = '  +2.32%  '
$ws.Range('D31').Value = '4.379'
$ws.Range('E31').Value = '  -0.80%  '
$ws.Range('E32').Value = '  +0.87%  '
$ws.Range('D33').Value = '4.063'
$ws.Range('E33').Value = '  -1.22%  '
$ws.Range('D34').Value = '0.04749'
$ws.Range('E34').Value = '  +1.91%  '
$ws.Range('D35').Value = '1.140'
$ws.Range('E35').Value = '  +4.15%  '
$ws.Range('D36').Value = '0.7047'
$ws.Range('E36').Value = '  +0.36%  '
$ws.Range('D37').Value = '2.707'
$ws.Range('E37').Value = '  +0.09%  '
$ws.Range('D38').Value = '0.01877'
$ws.Range('E38').Value = '  +1.27%  '
$ws.Range('D39').Value = '2.620'
$ws.Range('E39').Value = '  +4.24%  '
$ws.Range('D40').Value = '76.61'
$ws.Range('E40').Value = '  +4.14%  '
$ws.Range('D41').Value = '6.298'
$ws.Range('E41').Value = '  +0.27%  '
$ws.Range('D42').Value = '1.955'
$ws.Range('E42').Value = '  +1.70%  '
$ws.Range('D43').Value = '0.8505'
$ws.Range('E43').Value = '  -0.59%  '
$ws.Range('D44').Value = '0.4175'
$ws.Range('E44').Value = '  +0.37%  '
$ws.Range('D45').Value = '0.9998'
$ws.Range('E45').Value = '  +0.07%  '
$ws.Range('D46').Value = '103.12'
$ws.Range('E46').Value = '  -0.60%  '
$ws.Range('D47').Value = '7.206'
$ws.Range('E47').Value = '  +0.13%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '9.297'
$ws.Range('E48').Value = '  +0.49%  '
$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').Value = '947.90'
$ws.Range('E49').Value = '  -5.03%  '
$ws.Range('D50').Value = '34.23'
$ws.Range('E50').Value = '  +0.06%  '
$ws.Range('D51').Value = '0.05641'
$ws.Range('E51').Value = '  +0.13%  '
